# Experiment using universal quantifier in ms power point.
#
# Slide 34 ("Comparing two implementations of a function"):
#   Replace the "big wedge" (U+2A5D) operator with the universal
#   quantifier symbol "for all" (U+2200) in the subtitle text box.
#
# Slide 35 ("Comparing the same function ..."):
#   Re-type the words "that transforms" together with the preceding
#   space so they live in a single run.
#
# Slide 41 ("Challenges for Alabaster"):
#   Re-type "Could be of value for all kinds of mocking tools and
#   Monte Carlo studies." and " formulae, commuting diagrams . . .  "
#   so each lives in a single run instead of several split runs.
#
# Slide 42 ("Finally"):
#   Re-type "Transforms testing into development" as a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 34: swap the quantifier glyph
# ---------------------------------------------------------------
$s34 = $p.Slides.Item(34)
$shp34 = $s34.Shapes.Item(2)
$tr34 = $shp34.TextFrame.TextRange
$quant = $tr34.Characters(45, 2)
$quant.Text = [string]([char]0x2200) + " "

# ---------------------------------------------------------------
# Slide 35: merge " " + "that transforms " into one run
# ---------------------------------------------------------------
$s35 = $p.Slides.Item(35)
$shp35 = $s35.Shapes.Item(2)
$tr35 = $shp35.TextFrame.TextRange
$seg35 = $tr35.Characters(123, 17)
$seg35.Text = " that transforms "

# ---------------------------------------------------------------
# Slide 41: merge the three "Could "/"be of "/"value..." runs,
# and merge " formulae, commuting diagrams "/". . .  " runs
# ---------------------------------------------------------------
$s41 = $p.Slides.Item(41)
$shp41 = $s41.Shapes.Item(2)
$tr41 = $shp41.TextFrame.TextRange

$seg41a = $tr41.Characters(22, 73)
$seg41a.Text = "Could be of value for all kinds of mocking tools and Monte Carlo studies."

$seg41b = $tr41.Characters(254, 37)
$seg41b.Text = " formulae, commuting diagrams . . .  "

# ---------------------------------------------------------------
# Slide 42: merge "Transforms testing into "/"development" runs
# ---------------------------------------------------------------
$s42 = $p.Slides.Item(42)
$shp42 = $s42.Shapes.Item(2)
$tr42 = $shp42.TextFrame.TextRange
$seg42 = $tr42.Characters(19, 35)
$seg42.Text = "Transforms testing into development"
